$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N header (row 4): 2020
$ws.Cells.Item(4, 14).Value = 2020

# Data rows for column N (6-10)
$ws.Cells.Item(6, 14).Value = 1713
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(8, 14).Value = 379
$ws.Cells.Item(9, 14).Value = 180
$ws.Cells.Item(10, 14).Value = 798

# Selection
$ws.Range("L22").Select()
